$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Sex"), shifting Biomarker 1-5 from C:G to D:H
$ws.Columns("C:C").Insert()

# Header for the new column
$ws.Range("C1").Value = "Sex"

# Fill Sex values: even Index -> Male, odd Index -> Female (Index is in column A, rows 2-26)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    if (($idx % 2) -eq 0) {
        $ws.Cells.Item($r, 3).Value = "Male"
    } else {
        $ws.Cells.Item($r, 3).Value = "Female"
    }
}

# Update the selected cell to reflect the saved view state
$ws.Range("D1").Select()
